$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.176720142364502
$ws.Range("B1").Value = 2.41304349899292
$ws.Range("D1").Value = 2.337858915328979
$ws.Range("E1").Value = 1.200987696647644
